$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect so the cell contents can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (row 18, col A)
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$values = @{
    2  = @(0.05653192658765117, 0.01009443178117864)
    3  = @(0.02386177461388565, 0.006171648987463918)
    4  = @(0.03080959756799899, 0.009410409064720682)
    5  = @(0.03260243568267226, 0.009578544061302763)
    6  = @(0.03679558621270016, 0.002680102915951998)
    7  = @(0.01880409085290969, 0.004650024473813019)
    8  = @(0.004453496658891345, 0.001771479185119551)
    9  = @(0.006870245067229786, 0.00535885167464123)
    10 = @(0.07368584374498702, 0.002141327623126132)
    11 = @(0.07380418289447041, 0.001068947087119021)
    12 = @(0.1448155618611329, 0.003704510786664006)
    13 = @(0.382939570770756, 0.0009631380789771882)
    14 = @(0.1140256874847145, 0.005448601525608554)
    15 = @(0.9999999999999999, 0.003314153625254113)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# Restore sheet protection (matching the original protected state).
$ws.Protect()
